# Fig27.xlsx: 2017-02-13 snapshot - chunk 30
# Roll the "Short-Term Energy Outlook" title/source references from
# January 2017 to February 2017, and refresh the 2015-2018 forecast
# values for "Energy expenditures as share of GDP" (rows 51-54).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig27")

# Title cell (A2) and source-note cell (A55)
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("A55").Value = "Source: Short-Term Energy Outlook, February 2017."

# Updated forecast data points (2015-2018), column B
$ws.Range("B51").Value = 0.060564503439
$ws.Range("B52").Value = 0.054315399412
$ws.Range("B53").Value = 0.057963879846
$ws.Range("B54").Value = 0.057591287112
